{"js": "// Several of the phrases touched by this edit also occur, unmodified,\n// elsewhere in the document, so every search below is scoped to the one\n// paragraph that needs to change (found first by its distinctive text).\n\nasync function editParagraphContaining(marker, fn) {\n  const paragraphs = context.document.body.paragraphs;\n  paragraphs.load(\"text\");\n  await context.sync();\n  for (const p of paragraphs.items) {\n    if (p.text.indexOf(marker) !== -1) {\n      await fn(p);\n      return true;\n    }\n  }\n  return false;\n}\n\n// --- 1) \"Ch\u1ee9c n\u0103ng d\u00e0nh cho ch\u1ee7 nh\u00e0.\" -> \"Ch\u1ee9c n\u0103ng d\u00e0nh cho ng\u01b0\u1eddi cho thu\u00ea.\" ---\nawait editParagraphContaining(\"Ch\u1ee9c n\u0103ng d\u00e0nh cho ch\u1ee7 nh\u00e0\", async (p) => {\n  const results = p.search(\"ch\u1ee7 nh\u00e0\", { matchCase: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"ng\u01b0\u1eddi cho thu\u00ea\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n});\n\n// --- 2) \"Ch\u1ee9c n\u0103ng d\u00e0nh cho qu\u00e0n l\u00fd.\" -> \"Ch\u1ee9c n\u0103ng d\u00e0nh cho qu\u1ea3n l\u00fd.\" ---\nawait editParagraphContaining(\"Ch\u1ee9c n\u0103ng d\u00e0nh cho qu\u00e0n l\u00fd\", async (p) => {\n  const results = p.search(\"qu\u00e0n l\u00fd\", { matchCase: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"qu\u1ea3n l\u00fd\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n});\n\n// --- 3) Remove \" d\u01b0\u1edbi 3 gi\u00e2y cho c\u00e1c trang ch\u00ednh,\" from the page-speed bullet ---\nawait editParagraphContaining(\"T\u1ed1c \u0111\u1ed9 load trang nhanh\", async (p) => {\n  const results = p.search(\" d\u01b0\u1edbi 3 gi\u00e2y cho c\u00e1c trang ch\u00ednh,\", { matchCase: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n});\n\n// --- 4) \"Google: \u0110\u0103ng nh\u1eadp website b\u1eb1ng t\u00e0i kho\u1ea3n google\" runs get merged ---\n// Visible text is unchanged; rewrite the paragraph so it collapses back\n// into a single run (matching the canonical OOXML after the edit).\nawait editParagraphContaining(\"Google: \u0110\u0103ng nh\u1eadp website b\u1eb1ng t\u00e0i kho\u1ea3n google\", async (p) => {\n  p.load(\"text\");\n  await context.sync();\n  p.insertText(p.text, Word.InsertLocation.replace);\n  await context.sync();\n});\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Get-ParagraphRangeByMarker($doc, [string]$marker) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $para = $doc.Paragraphs.Item($i)\n        if ($para.Range.Text.Contains($marker)) {\n            return $para.Range\n        }\n    }\n    return $null\n}\n\nfunction Replace-InRange($range, [string]$findText, [string]$replaceText) {\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0            # wdFindStop - stay inside the supplied range\n    $find.Format = $false\n    $find.MatchCase = $false\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n\n# --- 1) \"Ch\u1ee9c n\u0103ng d\u00e0nh cho ch\u1ee7 nh\u00e0.\" -> \"Ch\u1ee9c n\u0103ng d\u00e0nh cho ng\u01b0\u1eddi cho thu\u00ea.\" ---\n# (several other paragraphs also contain \"ch\u1ee7 nh\u00e0\" unchanged, so we scope the\n# find/replace to just this heading paragraph)\n$range1 = Get-ParagraphRangeByMarker $d \"Ch\u1ee9c n\u0103ng d\u00e0nh cho ch\u1ee7 nh\u00e0\"\nif ($range1 -ne $null) {\n    Replace-InRange $range1 \"ch\u1ee7 nh\u00e0\" \"ng\u01b0\u1eddi cho thu\u00ea\"\n}\n\n# --- 2) \"Ch\u1ee9c n\u0103ng d\u00e0nh cho qu\u00e0n l\u00fd.\" -> \"Ch\u1ee9c n\u0103ng d\u00e0nh cho qu\u1ea3n l\u00fd.\" ---\n$range2 = Get-ParagraphRangeByMarker $d \"Ch\u1ee9c n\u0103ng d\u00e0nh cho qu\u00e0n l\u00fd\"\nif ($range2 -ne $null) {\n    Replace-InRange $range2 \"qu\u00e0n l\u00fd\" \"qu\u1ea3n l\u00fd\"\n}\n\n# --- 3) Remove \" d\u01b0\u1edbi 3 gi\u00e2y cho c\u00e1c trang ch\u00ednh,\" from the page-speed bullet ---\n$range3 = Get-ParagraphRangeByMarker $d \"T\u1ed1c \u0111\u1ed9 load trang nhanh\"\nif ($range3 -ne $null) {\n    Replace-InRange $range3 \" d\u01b0\u1edbi 3 gi\u00e2y cho c\u00e1c trang ch\u00ednh,\" \"\"\n}\n\n# --- 4) \"Google: \u0110\u0103ng nh\u1eadp website b\u1eb1ng t\u00e0i kho\u1ea3n google\" runs get merged ---\n# Visible text is unchanged; re-running Find/Replace with the identical text\n# collapses the paragraph back down to a single run (matching the canonical\n# OOXML after the edit).\n$range4 = Get-ParagraphRangeByMarker $d \"Google: \u0110\u0103ng nh\u1eadp website b\u1eb1ng t\u00e0i kho\u1ea3n google\"\nif ($range4 -ne $null) {\n    Replace-InRange $range4 \"Google: \u0110\u0103ng nh\u1eadp website b\u1eb1ng t\u00e0i kho\u1ea3n google\" \"Google: \u0110\u0103ng nh\u1eadp website b\u1eb1ng t\u00e0i kho\u1ea3n google\"\n}\n"}
